# Update numeric "sold/viewed" counters (column F) and one ticket-price
# value (column G2) across the four worksheets, matching the regenerated
# gh-pages data snapshot at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# ---- Sheet 1: 展览 ----
$ws = $wb.Worksheets.Item(1)
$ws.Range("G2").Value = 88
$ws.Range("F3").Value = 967
$ws.Range("F5").Value = 8607
$ws.Range("F6").Value = 164
$ws.Range("F7").Value = 49
$ws.Range("F8").Value = 1947
$ws.Range("F9").Value = 6080
$ws.Range("F10").Value = 602
$ws.Range("F13").Value = 8639
$ws.Range("F14").Value = 10227
$ws.Range("F15").Value = 1182
$ws.Range("F16").Value = 1027
$ws.Range("F17").Value = 4739
$ws.Range("F18").Value = 749
$ws.Range("F19").Value = 368
$ws.Range("F21").Value = 309
$ws.Range("F22").Value = 170
$ws.Range("F23").Value = 1274
$ws.Range("F24").Value = 177
$ws.Range("F25").Value = 1836
$ws.Range("F26").Value = 815
$ws.Range("F27").Value = 1117
$ws.Range("F28").Value = 801
$ws.Range("F29").Value = 1966
$ws.Range("F30").Value = 383
$ws.Range("F31").Value = 547
$ws.Range("F32").Value = 2501
$ws.Range("F33").Value = 304
$ws.Range("F34").Value = 145
$ws.Range("F35").Value = 1599
$ws.Range("F38").Value = 38
$ws.Range("F39").Value = 863
$ws.Range("F40").Value = 554
$ws.Range("F41").Value = 3160
$ws.Range("F44").Value = 469
$ws.Range("F45").Value = 552
$ws.Range("F47").Value = 879
$ws.Range("F48").Value = 215
$ws.Range("F49").Value = 4157

# ---- Sheet 2: 演出 ----
$ws = $wb.Worksheets.Item(2)
$ws.Range("F18").Value = 5

# ---- Sheet 3: 本地生活 ----
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 5609

# ---- Sheet 4: 全部类型 ----
$ws = $wb.Worksheets.Item(4)
$ws.Range("G2").Value = 88
$ws.Range("F3").Value = 967
$ws.Range("F5").Value = 8607
$ws.Range("F6").Value = 164
$ws.Range("F7").Value = 49
$ws.Range("F9").Value = 6080
$ws.Range("F10").Value = 602
$ws.Range("F11").Value = 8639
$ws.Range("F12").Value = 10227
$ws.Range("F14").Value = 1182
$ws.Range("F15").Value = 1027
$ws.Range("F16").Value = 4739
$ws.Range("F17").Value = 749
$ws.Range("F18").Value = 368
$ws.Range("F20").Value = 309
$ws.Range("F22").Value = 170
$ws.Range("F23").Value = 1274
$ws.Range("F24").Value = 177
$ws.Range("F25").Value = 1836
$ws.Range("F26").Value = 815
$ws.Range("F27").Value = 802
$ws.Range("F29").Value = 1966
$ws.Range("F30").Value = 383
$ws.Range("F31").Value = 547
$ws.Range("F32").Value = 2501
$ws.Range("F36").Value = 5
$ws.Range("F38").Value = 863
$ws.Range("F40").Value = 554
$ws.Range("F44").Value = 469
$ws.Range("F45").Value = 552
$ws.Range("F46").Value = 879
$ws.Range("F47").Value = 215
$ws.Range("F48").Value = 4157
